$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update Runmode column: stop running A Suite, start running E and F suites
$ws.Range("C2").Value = "N"
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"

# Update selection to C8
$ws.Range("C8").Select()
